$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new columns are inserted right after the (blank) label column:
#   B -> "$ bold(All)"   (new)
#   C -> "$ bold(Europe)" (new)
# The existing country columns shift right and get reordered, and the old
# "United States" header/column becomes "USA".
$ws.Range("B1").Value = '$ bold(All)'
$ws.Range("C1").Value = '$ bold(Europe)'
$ws.Range("D1").Value = "France"
$ws.Range("E1").Value = "Germany"
$ws.Range("F1").Value = "Italy"
$ws.Range("G1").Value = "Poland"
$ws.Range("H1").Value = "Spain"
$ws.Range("I1").Value = "United Kingdom"
$ws.Range("J1").Value = "Switzerland"
$ws.Range("K1").Value = "Japan"
$ws.Range("L1").Value = "Saudi Arabia"
$ws.Range("M1").Value = "USA"

# Data row 2 with the new values (including the two new columns).
$ws.Range("A2").Value = "Global climate scheme (GCS)"
$ws.Range("B2").Value = 0.55468268979227
$ws.Range("C2").Value = 0.622017379754577
$ws.Range("D2").Value = 0.612356985514925
$ws.Range("E2").Value = 0.581777133133904
$ws.Range("F2").Value = 0.739913584563118
$ws.Range("G2").Value = 0.534948999740369
$ws.Range("H2").Value = 0.664453305286869
$ws.Range("I2").Value = 0.592086655549235
$ws.Range("J2").Value = 0.638546618131183
$ws.Range("K2").Value = 0.55247423355968
$ws.Range("L2").Value = 0.842896551734276
$ws.Range("M2").Value = 0.453059461603854
